$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("List1")

$ws.Range("A2").Value = "Zakaznik_1_Rokytná"
$ws.Range("A3").Value = "Zakaznik_2_Krumlov"
$ws.Range("A4").Value = "Zakaznik_3_Polánka"
$ws.Range("A9").Value = "Zakaznik_8_Krumlov_2"
$ws.Range("A8").Value = "Zakaznik_7_Budkovice"
$ws.Range("A7").Value = "Zakaznik_6_Letkovice"
$ws.Range("A6").Value = "Zakaznik_5_Ivančice"
$ws.Range("A5").Value = "Zakaznik_4_Dobřínsko"
$ws.Range("A11").Value = "Zakaznik_10_Vémyslice"
$ws.Range("A10").Value = "Zakaznik_9_Dobelice"

$ws.Range("A11").Select()
